# Update the odds values on row 2 (F2:AO2) of the active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = [ordered]@{
    "F2"  = 2.66
    "G2"  = 2.68
    "H2"  = 3.3
    "I2"  = 3.35
    "J2"  = 3.05
    "K2"  = 3.1
    "L2"  = 0
    "M2"  = 0
    "N2"  = 4.7
    "O2"  = 1.26
    "P2"  = 1.9
    "Q2"  = 2.06
    "R2"  = 1.25
    "S2"  = 4.8
    "T2"  = 1.29
    "U2"  = 3.95
    "V2"  = 1.42
    "W2"  = 1.58
    "X2"  = 1000
    "Y2"  = 5
    "Z2"  = 34
    "AA2" = 390
    "AB2" = 1000
    "AC2" = 3.65
    "AD2" = 21
    "AE2" = 270
    "AF2" = 1000
    "AG2" = 4.6
    "AH2" = 25
    "AI2" = 1000
    "AJ2" = 1000
    "AK2" = 12.5
    "AL2" = 55
    "AM2" = 1000
    "AN2" = 26
    "AO2" = 1000
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
